$wb = $excel.ActiveWorkbook

# --- Update the "Paths" sheet content ---
$ws = $wb.Worksheets.Item("Paths")

# Row 3: new description for the data-folder row, plus a value and variable name
$ws.Range("B3").Value = "Ordner für die einzulesenden Daten"
$ws.Range("C3").Value = "Daten_Input"
$ws.Range("D3").Value = "r_DIR"

# Rows 4-6: drop the "Daten_Input\" folder prefix from the file paths
# (now that the folder itself is captured separately in C3/r_DIR)
$ws.Range("C4").Value = "Beispiel_Lastgang_einlesen.xlsx"
$ws.Range("C5").Value = "Beispiel_PV_Input_aus_Polysun.xlsx"
$ws.Range("C6").Value = "Beispiel_LKW_Fahrdaten.xlsx"

# Selection on the Paths sheet moves to C5
$ws.Range("C5").Select()

# Activate the "Paths" sheet (it becomes the active tab)
$ws.Activate()
